$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "marks" column (column G) entirely - header + per-question values.
# Everything to its right (imageURL column, the instruction box, etc.) shifts
# one column to the left.
$ws.Columns("G").Delete()

# The hyperlinks' anchor cells don't auto-shift with the column delete, so
# recreate them pointing at their new (post-shift) locations.
$target1 = $ws.Range("K8")
$target2 = $ws.Range("G9")
$rId1Url = "https://safe-exam-admin-nuv.vercel.app/uploadImage.html"
$rId2Url = "https://i.ibb.co/jk1XdnBH/img4.jpg"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($target1, $rId1Url)
$ws.Hyperlinks.Add($target2, $rId2Url)

# Adding a hyperlink resets the cell to a generic hyperlink font; restore the
# original (pre-existing) look of these two cells so their formatting stays
# exactly as it was before the column shift.
$target1.Font.Size = 11
$target2.Font.Size = 10

# Match the saved selection state (whole column G selected, cursor at G1).
$ws.Columns("G").Select() | Out-Null
